$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '35.171.04'
Set-TextValue 'E2' '  -0.18%  '
Set-TextValue 'D3' '1.904.36'
Set-TextValue 'E3' '  +0.39%  '
Set-TextValue 'E4' '  -0.42%  '
Set-TextValue 'D5' '254.12'
Set-TextValue 'E5' '  +3.09%  '
Set-TextValue 'D6' '0.698'
Set-TextValue 'E6' '  +1.79%  '
Set-TextValue 'E7' '  -0.35%  '
Set-TextValue 'D8' '41.81'
Set-TextValue 'E8' '  +2.65%  '
Set-TextValue 'E9' '  +2.18%  '
Set-TextValue 'D10' '52.64'
Set-TextValue 'E10' '  +1.24%  '
Set-TextValue 'D11' '0.0758'
Set-TextValue 'E11' '  +5.33%  '
Set-TextValue 'E12' '  -0.47%  '
Set-TextValue 'D13' '13.28'
Set-TextValue 'E13' '  +4.96%  '
Set-TextValue 'D14' '2.181.18'
Set-TextValue 'E14' '  +0.43%  '
Set-TextValue 'D15' '0.737'
Set-TextValue 'E15' '  +3.86%  '
Set-TextValue 'E16' '  +4.84%  '
Set-TextValue 'D17' '1.901.20'
Set-TextValue 'E17' '  +0.17%  '
Set-TextValue 'D18' '35.175.50'
Set-TextValue 'E18' '  -0.15%  '
Set-TextValue 'D19' '73.85'
Set-TextValue 'E19' '  +2.62%  '
Set-TextValue 'D20' '0.0₃0843'
Set-TextValue 'E20' '  +2.98%  '
Set-TextValue 'D21' '243.41'
Set-TextValue 'E21' '  +1.19%  '
Set-TextValue 'D22' '13.07'
Set-TextValue 'E22' '  +2.72%  '
Set-TextValue 'D23' '5.05'
Set-TextValue 'E23' '  +5.59%  '
Set-TextValue 'E24' '  -0.38%  '
Set-TextValue 'D25' '2.44'
Set-TextValue 'E25' '  +5.04%  '
Set-TextValue 'D26' '2.33'
Set-TextValue 'E26' '  -1.13%  '
Set-TextValue 'D27' '167.98'
Set-TextValue 'E27' '  +0.18%  '
Set-TextValue 'E28' '  +0.14%  '
Set-TextValue 'D29' '18.55'
Set-TextValue 'E29' '  +1.24%  '
Set-TextValue 'E30' '  +0.09%  '
Set-TextValue 'D31' '4.128.23'
Set-TextValue 'E32' '  +6.76%  '
Set-TextValue 'D33' '4.34'
Set-TextValue 'E33' '  +4.49%  '
Set-TextValue 'D34' '0.0599'
Set-TextValue 'E34' '  +5.68%  '
Set-TextValue 'E35' '  +9.69%  '
Set-TextValue 'D36' '4.26'
Set-TextValue 'E36' '  +3.66%  '
Set-TextValue 'E37' '  -0.43%  '
Set-TextValue 'D38' '0.854'
Set-TextValue 'E38' '  -6.48%  '
Set-TextValue 'E39' '  -0.08%  '
Set-TextValue 'B40' 'InjectiveProtocol'
Set-TextValue 'C40' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D40' '17.29'
Set-TextValue 'E40' '  +5.30%  '
Set-TextValue 'B41' 'Aave'
Set-TextValue 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '98.38'
Set-TextValue 'E41' '  +8.13%  '
Set-TextValue 'E42' '  +4.28%  '
Set-TextValue 'E43' '  +1.87%  '
Set-TextValue 'E44' '  +1.36%  '
Set-TextValue 'D45' '2.43'
Set-TextValue 'E45' '  +0.05%  '
Set-TextValue 'D46' '1.305.73'
Set-TextValue 'E46' '  -3.07%  '
Set-TextValue 'E47' '  +0.13%  '
Set-TextValue 'B48' 'Gas'
Set-TextValue 'C48' 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
Set-TextValue 'D48' '12.36'
Set-TextValue 'E48' '  +1.75%  '
Set-TextValue 'B49' 'MXToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D49' '2.75'
Set-TextValue 'E49' '  -1.20%  '
Set-TextValue 'E50' '  +2.10%  '
Set-TextValue 'E51' '  +7.21%  '
